$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.264.33'
$ws.Range('E2').Value = '  +2.31%  '
$ws.Range('D3').Value = '3.386.74'
$ws.Range('E3').Value = '  +1.74%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.46'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.82%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '180.71'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.78%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +0.77%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.198'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +9.15%  '
$ws.Range('E10').Value = '  +1.46%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '48.84'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.04%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000285'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +5.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '685.73'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.18%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.62'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.13%  '
$ws.Range('D15').Value = '3.934.09'
$ws.Range('E15').Value = '  +1.70%  '
$ws.Range('D16').Value = '69.295.16'
$ws.Range('E16').Value = '  +2.43%  '
$ws.Range('E17').Value = '  +1.92%  '
$ws.Range('D18').Value = '3.382.14'
$ws.Range('E18').Value = '  +1.97%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.72'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.95%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.37'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.96%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.902'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.88%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.43'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.42%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.10'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.71%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '103.96'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.49%  '
$ws.Range('E25').Value = '  +1.55%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.73'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.59%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.62'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.47%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '34.32'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.87%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.70'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.71%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.98'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.95%  '
$ws.Range('B31').Value = 'dogwifhat'
$ws.Range('C31').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.70'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +11.74%  '
$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.19'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.93%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '555.39'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.21%  '
$ws.Range('E34').Value = '  +1.24%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '58.14'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.64%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.16%  '
$ws.Range('D37').Value = '3.701.26'
$ws.Range('E37').Value = '  -0.29%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.141'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +7.68%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '34.98'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.23%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.23'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.13%  '
$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.69'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.09%  '
$ws.Range('B42').Value = 'PEPE'
$ws.Range('C42').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D42').Value = '0.0₃0704'
$ws.Range('E42').Value = '  +4.47%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.339'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.17%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0421'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.03%  '
$ws.Range('E45').Value = '  -1.91%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.65'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.01%  '
$ws.Range('E47').Value = '  +0.78%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.39'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.06%  '
$ws.Range('B49').Value = 'FirstDigitalUSD'
$ws.Range('C49').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.00'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.14%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.50'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.43%  '
$ws.Range('E51').Value = '  -1.99%  '
